# Product Backlog.xlsx - "actualizar backlog e HT"
# Inserts a new backlog row (HT-0004 "Diagrama de actividades") above the
# existing row 8, pushing the rest of the table down by one row, widens
# column I to fit the new (longer) acceptance-criteria text, and updates
# the saved view (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new row at 8, shifting rows 8:49 down to 9:50 -----------
$ws.Rows.Item(8).Insert()

# Pull the formatting (style "s", borders, alignment, etc.) for the new
# row from the row right below it (which now holds the content that used
# to live in the original row 8), then stamp in the new HT-0004 data.
$ws.Range("B9:I9").Copy($ws.Range("B8:I8"))
$ws.Rows.Item(8).RowHeight = 51.6

$ws.Range("B8").Value = "HT-0004"
$ws.Range("C8").Value = "Como un desarrollador, necesito elaborar los diagramas de actividades de las funcionalidades del proyecto, con la finalidad de comprender con mayor claridad las condiciones y limitaciones asociadas con cada operación específica."
$ws.Range("D8").Value = "Diagrama de actividades"
$ws.Range("E8").Value = "Hecho"
$ws.Range("F8").Value = "3 puntos"
$ws.Range("G8").Value = "Sprint 3"
$ws.Range("H8").Value = "Alta"
$ws.Range("I8").Value = "Un diagrama de actividades por cada funcionalidad"

# --- Column I needs to be wider to hold the new acceptance-criteria text --
$ws.Columns.Item(9).ColumnWidth = 21

# --- Refresh the saved view (zoom level + selected cell) ---------------
$ws.Application.ActiveWindow.Zoom = 65
$ws.Range("L9").Select()
